$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "A 34759-2023" record (previously at row 12) is moved up to row 10 with
# updated signal-species data, pushing the old rows 10-11 down by one each.
# Net row count is unchanged, so we insert a row at 10 and then delete the
# (now duplicate, shifted to row 13) original "A 34759-2023" row.
$ws.Rows(10).Insert()
$ws.Rows(13).Delete()

# Populate the new row 10 with the updated "A 34759-2023" record: an extra
# signal species ("Dropptaggsvamp") was found, bumping Signalarter 3->4 and
# Alla arter 8->9.
$ws.Range("A10").Value = "A 34759-2023"
$ws.Range("B10").Value = 45139
$ws.Range("C10").Value = 45188
$ws.Range("D10").Value = "VÄSTMANLANDS LÄN"
$ws.Range("E10").Value = "SURAHAMMAR"
$ws.Range("F10").Value = "Bergvik skog väst AB"
$ws.Range("G10").Value = 38.3
$ws.Range("H10").Value = 3
$ws.Range("I10").Value = 4
$ws.Range("J10").Value = 3
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 3
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 9
$ws.Range("R10").Value = "Grantaggsvamp`r`nMotaggsvamp`r`nSkogshare`r`nDropptaggsvamp`r`nGrönpyrola`r`nMindre märgborre`r`nPlattlummer`r`nMattlummer`r`nRevlummer"
$ws.Range("S10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/artfynd/A 34759-2023.xlsx", "A 34759-2023")'
$ws.Range("T10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/kartor/A 34759-2023.png", "A 34759-2023")'
$ws.Range("V10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/klagomål/A 34759-2023.docx", "A 34759-2023")'
$ws.Range("W10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/klagomålsmail/A 34759-2023.docx", "A 34759-2023")'
$ws.Range("X10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/tillsyn/A 34759-2023.docx", "A 34759-2023")'
$ws.Range("Y10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/tillsynsmail/A 34759-2023.docx", "A 34759-2023")'

# Row 10 keeps the standard 15pt row height (Insert() can otherwise autofit to
# the wrapped multi-line text we just wrote into R10)
$ws.Rows(10).RowHeight = 15

# Every data row's "Förändrad" (column C) timestamp moves from 45186 to 45188
for ($i = 2; $i -le 199; $i++) {
    $ws.Cells.Item($i, 3).Value = 45188
}
